$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the special "final row" number format (date-only) that
# currently lives on A13, before we touch anything.
$lastRowFormat = $ws.Range("A13").NumberFormat

# Row 13 reverts to the regular date+time format shared by every other
# data row (copy it from a sibling row, e.g. A12).
$ws.Range("A13").NumberFormat = $ws.Range("A12").NumberFormat

# Append the new day's data as row 14.
$ws.Range("A14").Value = 45963
$ws.Range("B14").Value = 26
$ws.Range("C14").Value = 36
$ws.Range("D14").Value = 32

# Row 14 becomes the new "final row" and inherits the date-only format
# that row 13 previously had.
$ws.Range("A14").NumberFormat = $lastRowFormat
